# "exam b attempt 2"
# Add a new worksheet ("Exam B 2") at the end of the workbook, holding a
# fresh (re-shuffled) attempt at the "Exam B" question bank: the matching
# headers in columns D/F/I/K/N plus the user's picked answers in column A
# for the 85 questions. No answer-key / grading columns (B/C/D6) yet since
# this attempt hasn't been checked.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Add($null, $wb.Worksheets($wb.Worksheets.Count))
$ws.Name = "Exam B 2"

# --- Header / matching-question blocks -------------------------------
# Written column-by-column (D, then F, then I, then K, then N) so that the
# brand-new shared strings this introduces get interned in the same order
# they were in the source workbook.

$ws.Range("D1").Value = "2 No Yes No"
$ws.Range("D2").Value = "2 Yes No No"
$ws.Range("D3").Value = "2 Yes No Yes"

$ws.Range("F1").Value = "Remote wipe"
$ws.Range("F2").Value = "Face recognition"
$ws.Range("F3").Value = "Full device encryption"
$ws.Range("F4").Value = "Host-based firewall"
$ws.Range("F5").Value = "Anti-malware"
$ws.Range("F6").Value = "Smart card"

$ws.Range("I1").Value = "CPU registers"
$ws.Range("I2").Value = "Routing table"
$ws.Range("I3").Value = "Temporary files"
$ws.Range("I4").Value = "Event logs"
$ws.Range("I5").Value = "Backup tapes"

$ws.Range("K1").Value = "Hashing"
$ws.Range("K2").Value = "Digital signature"
$ws.Range("K3").Value = "Encryption"
$ws.Range("K4").Value = "Key escrow"
$ws.Range("K5").Value = "Certificate authority"
$ws.Range("K6").Value = "Perfect forward secrecy"

$ws.Range("N2").Value = "Data in transit"
$ws.Range("N1").Value = "Data at rest"
$ws.Range("N3").Value = "Data at rest"
$ws.Range("N4").Value = "Data in use"
$ws.Range("N5").Value = "Data at rest"
$ws.Range("N6").Value = "Data in transit"
$ws.Range("N7").Value = "Data in use"
$ws.Range("N8").Value = "Data in use"
$ws.Range("N9").Value = "Data at rest"
$ws.Range("N10").Value = "Data at rest"

# --- Column A: the 85 answers picked on this attempt ------------------

$answers = @(
    "A","D","A","B","D","C","D","D","C","C",
    "D","D","D","D","C","B","D","C","A","B",
    "C","C","A","B","C","D","A","A","C","B",
    "A","B","D","AF","B","C","B","A","AE","D",
    "D","C","D","B","C","C","C","B","B","E",
    "D","AB","C","D","C","B","D","A","B","D",
    "C","B","B","A","D","A","AD","A","B","A",
    "C","C","B","A","C","BE","A","D","A","A",
    "C","A","B","C","C"
)

for ($i = 0; $i -lt $answers.Length; $i++) {
    $ws.Cells.Item(6 + $i, 1).Value = $answers[$i]
}

# --- View state: make the new sheet active with the same selection/scroll
# position it had when the workbook was last saved ---------------------

$ws.Activate()
$ws.Range("B1").Select()
$ws.Range("K18").Select()
